# Update Fecha (D), Volumen (J), Precio minimo (K), Precio maximo (L),
# Precio promedio ponderado (M) and Precio $/Kg (P) for rows 2-16.
# Values below reflect the new (post-edit) contents for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2;  D = 44659; J = 25; K = 10000; L = 10000; M = 10000; P = 1000 },
    @{ Row = 3;  D = 44425; J = 30; K = 13000; L = 13000; M = 13000; P = 1300 },
    @{ Row = 4;  D = 44463; J = 25; K = 12000; L = 12000; M = 12000; P = 1200 },
    @{ Row = 5;  D = 44348; J = 20; K = 10000; L = 10000; M = 10000; P = 1000 },
    @{ Row = 6;  D = 44369; J = 25; K = 8000;  L = 8000;  M = 8000;  P = 800  },
    @{ Row = 7;  D = 44523; J = 30; K = 9000;  L = 9000;  M = 9000;  P = 900  },
    @{ Row = 8;  D = 44525; J = 20; K = 9000;  L = 9000;  M = 9000;  P = 900  },
    @{ Row = 9;  D = 44656; J = 25; K = 10000; L = 10000; M = 10000; P = 1000 },
    @{ Row = 10; D = 44526; J = 25; K = 9000;  L = 9000;  M = 9000;  P = 900  },
    @{ Row = 11; D = 44663; J = 30; K = 12000; L = 12000; M = 12000; P = 1200 },
    @{ Row = 12; D = 44469; J = 20; K = 12000; L = 12000; M = 12000; P = 1200 },
    @{ Row = 13; D = 44649; J = 25; K = 10000; L = 10000; M = 10000; P = 1000 },
    @{ Row = 14; D = 44530; J = 30; K = 10000; L = 10000; M = 10000; P = 1000 },
    @{ Row = 15; D = 44473; J = 25; K = 11000; L = 11000; M = 11000; P = 1100 },
    @{ Row = 16; D = 44645; J = 25; K = 10000; L = 10000; M = 10000; P = 1000 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value  = $r.D   # D - Fecha
    $ws.Cells.Item($r.Row, 10).Value = $r.J   # J - Volumen
    $ws.Cells.Item($r.Row, 11).Value = $r.K   # K - Precio minimo
    $ws.Cells.Item($r.Row, 12).Value = $r.L   # L - Precio maximo
    $ws.Cells.Item($r.Row, 13).Value = $r.M   # M - Precio promedio ponderado
    $ws.Cells.Item($r.Row, 16).Value = $r.P   # P - Precio $/Kg
}
